# "optimized jump and call opcodes"
#
# The 'Instruction Set' sheet lists, among other things, the conditional/
# unconditional JMP/JC/JNC/JZ/JNZ/CALL/CC/CNC/CZ/CNZ "given address" family
# of opcodes in column H (8-bit binary opcode). This commit re-assigns 20 of
# those opcodes (rows 110-129) to a denser encoding. Column I (hex opcode)
# is a formula (BIN2HEX(NUMBERVALUE(H),2)) so it recalculates automatically,
# and the 'Opcodes' sheet looks commands up by that hex value via
# INDEX/MATCH, so it also recalculates automatically once column H changes.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Instruction Set")

# Row -> new 8-bit opcode (stored as text, matching the existing "@" / Text
# number format already applied to column H on these rows).
$updates = [ordered]@{
    110 = "01110101"
    111 = "01110111"
    112 = "01111001"
    113 = "01111011"
    114 = "01111101"
    115 = "01111111"
    116 = "10001111"
    117 = "10011111"
    118 = "10101111"
    119 = "10111111"
    120 = "11000001"
    121 = "11000011"
    122 = "11000101"
    123 = "11000111"
    124 = "11001011"
    125 = "11001111"
    126 = "11010001"
    127 = "11010011"
    128 = "11010101"
    129 = "11010111"
}

foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 8).Value = $updates[$row]
}
